$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before existing row 1053; this shifts all data
# currently in rows 1053:1116 down to rows 1054:1117 (and grows the sheet's
# used range to A1:T1117), matching the committed diff.
$ws.Rows.Item(1053).Insert()

# Populate the freshly inserted row 1053 with the new daily price record
# (same market/product context as its neighbours, new date + grade/price).
$ws.Range("A1053").Value = 10
$ws.Range("B1053").Value = "Vega Modelo de Temuco"
$ws.Range("C1053").Value = "La Araucanía"
$ws.Range("D1053").Value = 44610
$ws.Range("E1053").Value = 9
$ws.Range("F1053").Value = "Fruta"
$ws.Range("G1053").Value = 100102
$ws.Range("H1053").Value = "Cítricos"
$ws.Range("I1053").Value = 100102003
$ws.Range("J1053").Value = "Limón"
$ws.Range("K1053").Value = "Sin especificar"
$ws.Range("L1053").Value = "2a amarillo"
$ws.Range("M1053").Value = 200
$ws.Range("N1053").Value = 15000
$ws.Range("O1053").Value = 15000
$ws.Range("P1053").Value = 15000
$ws.Range("Q1053").Value = "$/bandeja 15 kilos"
$ws.Range("R1053").Value = "Región de O'Higgins"
$ws.Range("S1053").Value = 1000
$ws.Range("T1053").Value = 15
